# Auto-generated edit script: updates the cryptos price/volume table
# to the values captured by the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / link / percentage updates (safe to assign directly; Excel
#     keeps these as text because they are not valid numeric literals) ---
$ws.Range("D2").Value = "25.925.21"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.635.57"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.861.85"
$ws.Range("D14").Value = "1.636.41"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "25.922.25"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "1.147.72"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "1.771.61"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E48").Value = "  +5.61%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E51").Value = "  +2.37%  "

# --- Price cells whose new text would otherwise be auto-parsed by Excel as a
#     number (e.g. "19.70" -> 19.7). Force text entry via a temporary "@"
#     number format, then clear the format again so the cell is left with the
#     same (default) style it started with, only the text value changes. ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.02"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.70"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.97"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.89"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.35"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.30"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.909"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0157"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.84"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.65"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0512"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.46"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.417"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.60"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0964"

$ws.Range("D5").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
